# working_hours.xlsx: bind/insert the new "24 Jul 2014, 22:00-23:00" record
# as row 152 (it was previously an empty separator row), pushing the
# separator + the three summary rows ("sum [min]", "sum [h]",
# "sum [working weeks]") down by one row, and refresh the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 152 - this shifts the old row 152 (blank
# separator) and the three summary rows below it down to 153-156,
# including their formulas/styles/shared strings.
$ws.Rows.Item(152).Insert()

# Fill the newly inserted row 152 with the new time-tracking entry.
$ws.Range("A152").Value = 2014
$ws.Range("B152").Value = 7
$ws.Range("C152").Value = 24
$ws.Range("D152").Value = 0.91666666666666663
$ws.Range("E152").Value = 0.95833333333333337
$ws.Range("F152").Formula = "=(E152-D152)*24*60"
$ws.Range("G152").Formula = "=F152/60"

# The summary formulas (now on rows 154-156) need their ranges/refs
# nudged to include the new data row 152.
$ws.Range("F154").Formula = "=SUM(F2:F152)"
$ws.Range("F155").Formula = "=F154/60"
$ws.Range("F156").Formula = "=F155/38.5"

# Match the author's final selection/cursor position.
$ws.Range("F152").Select() | Out-Null
